$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC30_TO_Verify_Pagi_sort_filter")

# Insert a new row above row 15, shifting existing rows 15-25 down to 16-26.
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(15).Insert()

# Populate the freshly inserted row 15 with the new step; clear the rest.
$ws.Range("B15").Value = "TINY_SCROLL_DOWN"
$ws.Range("C15").Value = $null
$ws.Range("D15").Value = $null
$ws.Range("E15").Value = $null

# Update the active selection to mirror the post-edit state.
$ws.Range("A7").Select()
$ws.Range("C27").Select()
